$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1204.6538
$ws.Range("I98").Value = 1243.8572
$ws.Range("J98").Value = 1040
$ws.Range("K98").Value = 1243.8572
$ws.Range("L98").Value = 1040
$ws.Range("M98").Value = 254.1428000000001
$ws.Range("N98").Value = -4036

$ws.Range("H106").Value = 3753.6365
$ws.Range("I106").Value = 3698.3333
$ws.Range("K106").Value = 3698.3333
$ws.Range("M106").Value = -3067.3333

$ws.Range("H122").Value = 1204.6538
$ws.Range("I122").Value = 1243.8572
$ws.Range("J122").Value = 1040
$ws.Range("K122").Value = 3731.5716
$ws.Range("L122").Value = 3120
$ws.Range("M122").Value = -1281.5716
$ws.Range("N122").Value = -8020

$ws.Range("H137").Value = 1251.3636
$ws.Range("I137").Value = 1013.5294
$ws.Range("J137").Value = 2060
$ws.Range("K137").Value = 3040.5882
$ws.Range("L137").Value = 6180
$ws.Range("M137").Value = -490.5882000000001
$ws.Range("N137").Value = -11280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5485027
$ws.Range("I32").Value = 6082615.5
$ws.Range("J32").Value = 7133.3335
$ws.Range("K32").Value = 6082615.5
$ws.Range("L32").Value = 7133.3335
$ws.Range("M32").Value = -6082328.5
$ws.Range("N32").Value = -7707.3335

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61:N61").ClearContents()

$ws.Range("H74").Value = 1213.1364
$ws.Range("I74").Value = 1354.3889
$ws.Range("K74").Value = 1354.3889
$ws.Range("M74").Value = -480.3888999999999

$ws.Range("H77").Value = 1213.1364
$ws.Range("I77").Value = 1354.3889
$ws.Range("K77").Value = 6771.9445
$ws.Range("M77").Value = -2403.9445

$ws.Range("H97").Value = 645.9375
$ws.Range("I97").Value = 587.3077
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 587.3077
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -91.30769999999995
$ws.Range("N97").Value = -1892

$ws.Range("H132").Value = 2149.7778
$ws.Range("I132").Value = 1666.5238
$ws.Range("J132").Value = 3841.1667
$ws.Range("K132").Value = 4999.5714
$ws.Range("L132").Value = 11523.5001
$ws.Range("M132").Value = -2469.5714
$ws.Range("N132").Value = -16583.5001

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136:N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 55590640
$ws.Range("I20").Value = 68192.78
$ws.Range("J20").Value = 111113090
$ws.Range("K20").Value = 68192.78
$ws.Range("L20").Value = 111113090
$ws.Range("M20").Value = -67945.78
$ws.Range("N20").Value = -111113584

$ws.Range("H64").Value = 778.1111
$ws.Range("I64").Value = 651
$ws.Range("J64").Value = 879.8
$ws.Range("K64").Value = 651
$ws.Range("L64").Value = 879.8
$ws.Range("M64").Value = -426
$ws.Range("N64").Value = -1329.8

$ws.Range("H67").Value = 778.1111
$ws.Range("I67").Value = 651
$ws.Range("J67").Value = 879.8
$ws.Range("K67").Value = 651
$ws.Range("L67").Value = 879.8
$ws.Range("M67").Value = 129
$ws.Range("N67").Value = -2439.8

$ws.Range("H94").Value = 1082.7
$ws.Range("I94").Value = 1057.2941
$ws.Range("J94").Value = 1226.6666
$ws.Range("K94").Value = 1057.2941
$ws.Range("L94").Value = 1226.6666
$ws.Range("M94").Value = -606.2941000000001
$ws.Range("N94").Value = -2128.6666

$ws.Range("H99").Value = 2055.6667
$ws.Range("I99").Value = 2153.4546
$ws.Range("J99").Value = 980
$ws.Range("K99").Value = 2153.4546
$ws.Range("L99").Value = 980
$ws.Range("M99").Value = -655.4546
$ws.Range("N99").Value = -3976

$ws.Range("H107").Value = 8248.263000000001
$ws.Range("I107").Value = 1328.7222
$ws.Range("J107").Value = 132800
$ws.Range("K107").Value = 1328.7222
$ws.Range("L107").Value = 132800
$ws.Range("M107").Value = 591.2778000000001
$ws.Range("N107").Value = -136640

$ws.Range("H134").Value = 30148.945
$ws.Range("I134").Value = 2353.44
$ws.Range("J134").Value = 93320.55
$ws.Range("K134").Value = 7060.32
$ws.Range("L134").Value = 279961.65
$ws.Range("M134").Value = -4525.32
$ws.Range("N134").Value = -285031.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 187.8
$ws.Range("I6").Value = 187.8
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 187.8
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -74.80000000000001
$ws.Range("N6").ClearContents()

$ws.Range("H31").Value = 2643.62
$ws.Range("J31").Value = 2630.6
$ws.Range("L31").Value = 2630.6
$ws.Range("N31").Value = -3220.6

$ws.Range("H34").Value = 2643.62
$ws.Range("J34").Value = 2630.6
$ws.Range("L34").Value = 2630.6
$ws.Range("N34").Value = -3034.6

$ws.Range("H58").Value = 4716.2593
$ws.Range("I58").Value = 932.93335
$ws.Range("J58").Value = 9445.416999999999
$ws.Range("K58").Value = 932.93335
$ws.Range("L58").Value = 9445.416999999999
$ws.Range("M58").Value = -729.93335
$ws.Range("N58").Value = -9851.416999999999

$ws.Range("H134").Value = 2944.375
$ws.Range("I134").Value = 555
$ws.Range("J134").Value = 3285.7144
$ws.Range("K134").Value = 1665
$ws.Range("L134").Value = 9857.143199999999
$ws.Range("M134").Value = 870
$ws.Range("N134").Value = -14927.1432

$ws.Range("H136").Value = 4716.2593
$ws.Range("I136").Value = 932.93335
$ws.Range("J136").Value = 9445.416999999999
$ws.Range("K136").Value = 2798.80005
$ws.Range("L136").Value = 28336.251
$ws.Range("M136").Value = -248.8000499999998
$ws.Range("N136").Value = -33436.251

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1350
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1350
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 4050
$ws.Range("N92").Value = -6546
$ws.Range("M92").ClearContents()

$ws.Range("H131").Value = 783.97
$ws.Range("J131").Value = 840.97754
$ws.Range("L131").Value = 2522.93262
$ws.Range("N131").Value = -12602.93262

$ws.Range("H137").Value = 51016.453
$ws.Range("J137").Value = 100266.55
$ws.Range("L137").Value = 300799.65
$ws.Range("M137").Value = -199.0907999999999
$ws.Range("N137").Value = -310999.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4177.75
$ws.Range("I70").Value = 4065.3572
$ws.Range("K70").Value = 4065.3572
$ws.Range("M70").Value = -3795.3572

$ws.Range("H73").Value = 4177.75
$ws.Range("I73").Value = 4065.3572
$ws.Range("K73").Value = 4065.3572
$ws.Range("M73").Value = -3129.3572

$ws.Range("H97").Value = 2801.1482
$ws.Range("I97").Value = 3209.4119
$ws.Range("J97").Value = 2107.1
$ws.Range("K97").Value = 3209.4119
$ws.Range("L97").Value = 2107.1
$ws.Range("M97").Value = -2713.4119
$ws.Range("N97").Value = -3099.1

$ws.Range("H126").Value = 2172.4614
$ws.Range("I126").Value = 1964.2222
$ws.Range("J126").Value = 2641
$ws.Range("K126").Value = 5892.6666
$ws.Range("L126").Value = 7923
$ws.Range("M126").Value = -3422.6666
$ws.Range("N126").Value = -12863

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5651.727
$ws.Range("I61").Value = 6402.5264
$ws.Range("J61").Value = 896.6667
$ws.Range("K61").Value = 6402.5264
$ws.Range("L61").Value = 896.6667
$ws.Range("M61").Value = -6200.5264
$ws.Range("N61").Value = -1300.6667

$ws.Range("H113").Value = 5651.727
$ws.Range("I113").Value = 6402.5264
$ws.Range("J113").Value = 896.6667
$ws.Range("K113").Value = 6402.5264
$ws.Range("L113").Value = 896.6667
$ws.Range("M113").Value = -4232.5264
$ws.Range("N113").Value = -5236.6667

$ws.Range("H136").Value = 1895.2
$ws.Range("I136").Value = 1318.7556
$ws.Range("J136").Value = 3624.5334
$ws.Range("K136").Value = 3956.2668
$ws.Range("L136").Value = 10873.6002
$ws.Range("M136").Value = -1406.2668
$ws.Range("N136").Value = -15973.6002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9642.857
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 11625
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 11625
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -12873

$ws.Range("H65").Value = 9642.857
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 11625
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 58125
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -64365

$ws.Range("H96").Value = 3512.5
$ws.Range("I96").Value = 3600
$ws.Range("J96").Value = 3250
$ws.Range("K96").Value = 3600
$ws.Range("L96").Value = 3250
$ws.Range("M96").Value = -2227
$ws.Range("N96").Value = -5996

$ws.Range("H122").Value = 637.9167
$ws.Range("I122").Value = 585
$ws.Range("J122").Value = 902.5
$ws.Range("K122").Value = 1755
$ws.Range("L122").Value = 2707.5
$ws.Range("M122").Value = 695
$ws.Range("N122").Value = -7607.5
